$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")
$ws.Activate()

# Mark Obadiah (row 32) and Jonah (row 33) as finished books ("Book done" = 1).
# The "Verses done" column (E) is a shared formula IF(F=1,C,0), so it will
# recalculate to the verse counts for those rows automatically.
$ws.Range("F32").Value = 1
$ws.Range("F33").Value = 1

# Scroll the view down to where work left off (Micah 2:8) and select the
# newly-finished cells, matching where the editor's cursor ended up.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E32:F33").Select()
